{"js": "// Fix the \"Demo\" hyperlink text for the Netflix-clone project. In the\n// underlying OOXML this hyperlink's visible text \"Demo\" was split across\n// three separate runs (\"De\" + \"m\" + \"o\"), which is what the \"download /\n// demo button\" bug report referred to. Re-typing it back into a single\n// contiguous \"Demo\" run (keeping the Hyperlink character style) is the fix.\n\nconst body = context.document.body;\n\n// There are three \"Demo\" hyperlinks in this resume (Netflix clone, Shoe\n// selling website, Weather App projects). Only the first one - right after\n// the Netflix-clone project description - is the broken/split one, so\n// disambiguate by looking at the paragraph that immediately precedes each\n// match.\nconst matches = body.search(\"Demo\", { matchCase: true });\nmatches.load(\"text\");\nawait context.sync();\n\nlet matchParagraphs = [];\nfor (let i = 0; i < matches.items.length; i++) {\n  const par = matches.items[i].paragraphs.getFirst();\n  par.load(\"text\");\n  matchParagraphs.push(par);\n}\nawait context.sync();\n\nlet prevParagraphs = [];\nfor (let i = 0; i < matchParagraphs.length; i++) {\n  const prev = matchParagraphs[i].getPrevious();\n  prev.load(\"text\");\n  prevParagraphs.push(prev);\n}\nawait context.sync();\n\nlet targetRange = null;\nfor (let i = 0; i < matches.items.length; i++) {\n  if (prevParagraphs[i].text.indexOf(\"Netflix clone\") !== -1) {\n    targetRange = matches.items[i];\n    break;\n  }\n}\n\nif (!targetRange) {\n  throw new Error(\"Could not find the Netflix-clone project's 'Demo' hyperlink.\");\n}\n\n// A plain \"Demo\" -> \"Demo\" replace is a no-op (identical text), which would\n// leave the three original runs untouched. Replacing with different text\n// first, then replacing that *whole* marker text with \"Demo\" right after,\n// forces the editor to actually rewrite the range - collapsing the three\n// split runs (\"De\"/\"m\"/\"o\") into a single run with the correct \"Demo\" text.\nconst marker = \"\\u00B6FIXUP\\u00B6\";\ntargetRange.insertText(marker, Word.InsertLocation.replace);\nawait context.sync();\n\nconst tmpMatches = context.document.body.search(marker, { matchCase: true });\ntmpMatches.load(\"text\");\nawait context.sync();\n\nif (tmpMatches.items.length === 0) {\n  throw new Error(\"Could not re-locate the 'Demo' hyperlink after the intermediate edit.\");\n}\n\ntmpMatches.items[0].insertText(\"Demo\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fix the \"Demo\" hyperlink text for the Netflix-clone project. In the\n# underlying OOXML this hyperlink's visible text \"Demo\" was split across\n# three separate runs (\"De\" + \"m\" + \"o\"), which is what the \"download /\n# demo button\" bug report referred to. Re-typing it back into a single\n# contiguous \"Demo\" run (keeping the Hyperlink character style) is the fix.\n\n$d = $word.ActiveDocument\n\n# There are three \"Demo\" hyperlinks in this resume (Netflix clone, Shoe\n# selling website, Weather App projects). Only the first one - right after\n# the Netflix-clone project description - is the broken/split one, so walk\n# every \"Demo\" match and disambiguate using the text that precedes it.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Demo\"\n$find.MatchCase = $true\n$find.Forward = $true\n\n$targetStart = -1\n$targetEnd = -1\n$found = $find.Execute()\n$guard = 0\nwhile ($found -and $guard -lt 50) {\n    $start = $rng.Start\n    $end = $rng.End\n    $precedingText = $d.Range(0, $start).Text\n    $tail = $precedingText.Substring([Math]::Max(0, $precedingText.Length - 300))\n    if ($tail.Contains(\"Netflix clone\")) {\n        $targetStart = $start\n        $targetEnd = $end\n    }\n    $rng.Collapse(0)\n    $found = $find.Execute()\n    $guard++\n}\n\nif ($targetStart -lt 0) {\n    throw \"Could not find the Netflix-clone project's 'Demo' hyperlink.\"\n}\n\n# A plain \"Demo\" -> \"Demo\" assignment is a no-op (identical text) and would\n# leave the three original runs untouched, so briefly change the text to\n# something different, then set the *whole* inserted range back to \"Demo\"\n# right after. Both edits land in a single rewritten run, collapsing the\n# three split runs (\"De\"/\"m\"/\"o\") into one correct \"Demo\" run.\n$fixupText = [char]0xB6 + \"FIXUP\" + [char]0xB6\n$target = $d.Range($targetStart, $targetEnd)\n$target.Text = $fixupText\n$target2 = $d.Range($targetStart, $targetStart + $fixupText.Length)\n$target2.Text = \"Demo\"\n"}
